# Scheduled-runner data refresh: update market price / profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of Shiva_Profits.xlsx.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4918
$ws.Range("I19").Value = 5541.5713
$ws.Range("J19").Value = 4372.375
$ws.Range("K19").Value = 5541.5713
$ws.Range("L19").Value = 4372.375
$ws.Range("M19").Value = -5366.5713
$ws.Range("N19").Value = -4722.375

$ws.Range("H53").Value = 781.9167
$ws.Range("J53").Value = 385.2857
$ws.Range("L53").Value = 385.2857
$ws.Range("N53").Value = -1659.2857

$ws.Range("H80").Value = 966.7857
$ws.Range("J80").Value = 1284.8334
$ws.Range("L80").Value = 3854.5002
$ws.Range("N80").Value = -5850.5002

$ws.Range("H83").Value = 966.7857
$ws.Range("J83").Value = 1284.8334
$ws.Range("L83").Value = 11563.5006
$ws.Range("N83").Value = -21547.5006

$ws.Range("H98").Value = 1470.8572
$ws.Range("I98").Value = 1053.8334
$ws.Range("K98").Value = 1053.8334
$ws.Range("M98").Value = 444.1666

$ws.Range("H116").Value = 4775.6816
$ws.Range("J116").Value = 4414.722
$ws.Range("L116").Value = 4414.722
$ws.Range("N116").Value = -11298.722

$ws.Range("H122").Value = 1470.8572
$ws.Range("I122").Value = 1053.8334
$ws.Range("K122").Value = 3161.5002
$ws.Range("M122").Value = -711.5001999999999

$ws.Range("H132").Value = 5797.685
$ws.Range("I132").Value = 3070.182
$ws.Range("K132").Value = 9210.545999999998
$ws.Range("M132").Value = -6680.545999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 49999
$ws.Range("J47").Value = 49999
$ws.Range("L47").Value = 49999
$ws.Range("N47").Value = -51449

$ws.Range("H54").Value = 40000
$ws.Range("I54").Value = 40000
$ws.Range("K54").Value = 40000
$ws.Range("M54").Value = -39231

$ws.Range("H61").Value = 5197.263
$ws.Range("I61").Value = 5861.871
$ws.Range("J61").Value = 2254
$ws.Range("K61").Value = 5861.871
$ws.Range("L61").Value = 2254
$ws.Range("M61").Value = -5649.871
$ws.Range("N61").Value = -2678

$ws.Range("H122").Value = 15644.083
$ws.Range("I122").Value = 1959
$ws.Range("J122").Value = 43014.25
$ws.Range("K122").Value = 5877
$ws.Range("L122").Value = 129042.75
$ws.Range("M122").Value = -3427
$ws.Range("N122").Value = -133942.75

$ws.Range("H136").Value = 5197.263
$ws.Range("I136").Value = 5861.871
$ws.Range("J136").Value = 2254
$ws.Range("K136").Value = 17585.613
$ws.Range("L136").Value = 6762
$ws.Range("M136").Value = -15035.613
$ws.Range("N136").Value = -11862

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1293.1
$ws.Range("I80").Value = 846.25
$ws.Range("J80").Value = 1591
$ws.Range("K80").Value = 846.25
$ws.Range("L80").Value = 1591
$ws.Range("M80").Value = 151.75
$ws.Range("N80").Value = -3587

$ws.Range("H83").Value = 1293.1
$ws.Range("I83").Value = 846.25
$ws.Range("J83").Value = 1591
$ws.Range("K83").Value = 4231.25
$ws.Range("L83").Value = 7955
$ws.Range("M83").Value = 760.75
$ws.Range("N83").Value = -17939

$ws.Range("H94").Value = 1606
$ws.Range("J94").Value = 1760
$ws.Range("L94").Value = 1760
$ws.Range("N94").Value = -2662

$ws.Range("H99").Value = 2097
$ws.Range("I99").Value = 2097
$ws.Range("K99").Value = 2097
$ws.Range("M99").Value = -599

$ws.Range("H134").Value = 4159.769
$ws.Range("I134").Value = 3706.2144
$ws.Range("K134").Value = 11118.6432
$ws.Range("M134").Value = -8583.643199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.636364
$ws.Range("J7").Value = 98.666664
$ws.Range("L7").Value = 98.666664
$ws.Range("N7").Value = -324.666664

$ws.Range("H31").Value = 1777.1052
$ws.Range("I31").Value = 1544.8928
$ws.Range("J31").Value = 2427.3
$ws.Range("K31").Value = 1544.8928
$ws.Range("L31").Value = 2427.3
$ws.Range("M31").Value = -1249.8928
$ws.Range("N31").Value = -3017.3

$ws.Range("H34").Value = 1777.1052
$ws.Range("I34").Value = 1544.8928
$ws.Range("J34").Value = 2427.3
$ws.Range("K34").Value = 1544.8928
$ws.Range("L34").Value = 2427.3
$ws.Range("M34").Value = -1342.8928
$ws.Range("N34").Value = -2831.3

$ws.Range("H58").Value = 1947.1305
$ws.Range("I58").Value = 1846.9048
$ws.Range("J58").Value = 2999.5
$ws.Range("K58").Value = 1846.9048
$ws.Range("L58").Value = 2999.5
$ws.Range("M58").Value = -1643.9048
$ws.Range("N58").Value = -3405.5

$ws.Range("H62").Value = 111114090
$ws.Range("I62").Value = 250002260
$ws.Range("J62").Value = 3559.4
$ws.Range("K62").Value = 250002260
$ws.Range("L62").Value = 3559.4
$ws.Range("M62").Value = -250001636
$ws.Range("N62").Value = -4807.4

$ws.Range("H65").Value = 111114090
$ws.Range("I65").Value = 250002260
$ws.Range("J65").Value = 3559.4
$ws.Range("K65").Value = 1250011300
$ws.Range("L65").Value = 17797
$ws.Range("M65").Value = -1250008180
$ws.Range("N65").Value = -24037

$ws.Range("H86").Value = 25759904
$ws.Range("I86").Value = 43591876
$ws.Range("K86").Value = 43591876
$ws.Range("M86").Value = -43590753

$ws.Range("H89").Value = 25759904
$ws.Range("I89").Value = 43591876
$ws.Range("K89").Value = 217959380
$ws.Range("M89").Value = -217953764

$ws.Range("H132").Value = 17829.924
$ws.Range("I132").Value = 7726.8335
$ws.Range("K132").Value = 23180.5005
$ws.Range("M132").Value = -20650.5005

$ws.Range("H135").Value = 156246.45
$ws.Range("J135").Value = 156246.45
$ws.Range("L135").Value = 156246.45
$ws.Range("N135").Value = -166386.45

$ws.Range("H136").Value = 1947.1305
$ws.Range("I136").Value = 1846.9048
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 5540.7144
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -2990.7144
$ws.Range("N136").Value = -14098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1150.125
$ws.Range("I13").Value = 81.8
$ws.Range("J13").Value = 2930.6667
$ws.Range("K13").Value = 245.4
$ws.Range("L13").Value = 8792.000100000001
$ws.Range("M13").Value = -77.39999999999998
$ws.Range("N13").Value = -9128.000100000001

$ws.Range("H18").Value = 784
$ws.Range("I18").Value = 784
$ws.Range("K18").Value = 2352
$ws.Range("M18").Value = -2183

$ws.Range("H108").Value = 1162.4
$ws.Range("I108").Value = 1162.4
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3487.2
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -607.2000000000003
$ws.Range("N108").Value = ""

$ws.Range("H121").Value = 1595.1111
$ws.Range("I121").Value = 319.4
$ws.Range("J121").Value = 2085.7693
$ws.Range("K121").Value = 958.1999999999999
$ws.Range("L121").Value = 6257.3079
$ws.Range("M121").Value = 351.8000000000001
$ws.Range("N121").Value = -8877.3079

$ws.Range("H139").Value = 3946.6667
$ws.Range("I139").Value = 1736.9231
$ws.Range("K139").Value = 5210.7693
$ws.Range("M139").Value = -70.76929999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 400.14285
$ws.Range("I107").Value = 530.8889
$ws.Range("J107").Value = 164.8
$ws.Range("K107").Value = 530.8889
$ws.Range("L107").Value = 164.8
$ws.Range("M107").Value = 1389.1111
$ws.Range("N107").Value = -4004.8

$ws.Range("H113").Value = 3981.8438
$ws.Range("I113").Value = 3930.9
$ws.Range("K113").Value = 3930.9
$ws.Range("M113").Value = -1760.9

$ws.Range("H126").Value = 4740.9375
$ws.Range("I126").Value = 4400.4644
$ws.Range("K126").Value = 13201.3932
$ws.Range("M126").Value = -10731.3932

$ws.Range("H132").Value = 6403.9546
$ws.Range("I132").Value = 4968.265
$ws.Range("J132").Value = 10542.117
$ws.Range("K132").Value = 14904.795
$ws.Range("L132").Value = 31626.351
$ws.Range("M132").Value = -12374.795
$ws.Range("N132").Value = -36686.351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 543.4828
$ws.Range("I55").Value = 399.875
$ws.Range("J55").Value = 720.2308
$ws.Range("K55").Value = 399.875
$ws.Range("L55").Value = 720.2308
$ws.Range("M55").Value = -226.875
$ws.Range("N55").Value = -1066.2308

$ws.Range("H127").Value = 53238.332
$ws.Range("J127").Value = 53238.332
$ws.Range("L127").Value = 53238.332
$ws.Range("N127").Value = -63158.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 22395.715
$ws.Range("J69").Value = 25678.334
$ws.Range("L69").Value = 25678.334
$ws.Range("N69").Value = -27176.334

$ws.Range("H72").Value = 22395.715
$ws.Range("J72").Value = 25678.334
$ws.Range("L72").Value = 77035.00199999999
$ws.Range("N72").Value = -84523.00199999999

Write-Output "Applied scheduled runner price/profit updates."